$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before) values for columns B..H across rows 2..9
$data = @{}
for ($r = 2; $r -le 9; $r++) {
    $row = @{
        B = $ws.Cells.Item($r, 2).Value2
        C = $ws.Cells.Item($r, 3).Value2
        D = $ws.Cells.Item($r, 4).Value2
        E = $ws.Cells.Item($r, 5).Value2
        G = $ws.Cells.Item($r, 7).Value2
        H = $ws.Cells.Item($r, 8).Value2
    }
    $data[$r] = $row
}

# Rows 2..9 correspond to index (A) 0..7.
# The data in columns B..H is cyclically shifted: the row whose index is i
# receives the B..H content that previously belonged to the row whose
# index is (i + 3) mod 8 (rows are 2-based, so row = index + 2).
for ($i = 0; $i -le 7; $i++) {
    $destRow = $i + 2
    $srcIndex = ($i + 3) % 8
    $srcRow = $srcIndex + 2
    $src = $data[$srcRow]

    $ws.Cells.Item($destRow, 2).Value2 = $src.B
    $ws.Cells.Item($destRow, 3).Value2 = $src.C
    $ws.Cells.Item($destRow, 4).Value2 = $src.D
    $ws.Cells.Item($destRow, 5).Value2 = $src.E
    $ws.Cells.Item($destRow, 7).Value2 = $src.G
    $ws.Cells.Item($destRow, 8).Value2 = $src.H
}
